$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row B (Right marking value)
$ws.Range("B11").Value = 5

# Update "Total" row B (Correct total) and E (Corr/total marks text)
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
